# Add a "2021" column (column R) to the 1.5.2 disaster-loss / GDP sheet,
# mirroring the layout of the existing year columns (D..Q = 2007..2020).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Formatting -------------------------------------------------------
# Column R should look exactly like column Q (same number formats,
# borders, fonts, etc.) for every row that currently has data in Q.
# Copy Q4:Q44 (header row through the last data row) and paste only the
# formatting into the new R4:R44 range.
$ws.Range("Q4:Q44").Copy()
$ws.Range("R4:R44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header -------------------------------------------------------------
$ws.Range("R4").Value = 2021

# --- Кыргыз Республикасы (Kyrgyz Republic) block ------------------------
$ws.Range("R6").Formula = "=R7/R8*100"
$ws.Range("R7").Value = 1931.83
$ws.Range("R8").Value = 739818.5

# --- Баткен облусу (Batken oblast) block --------------------------------
$ws.Range("R10").Formula = "=R11/R12*100"
$ws.Range("R11").Value = 1552.9
$ws.Range("R12").Value = 25048.6

# --- Жалал-Абад облусу (Djalal-Abad oblast) block -----------------------
$ws.Range("R14").Formula = "=R15/R16*100"
$ws.Range("R15").Value = 125.7
$ws.Range("R16").Value = 82213.899999999994

# --- Ысык-Көл облусу (Ysyk-Kul oblast) block -----------------------------
$ws.Range("R18").Formula = "=R19/R20*100"
$ws.Range("R19").Value = 99.6
$ws.Range("R20").Value = 80059.600000000006

# --- Нарын облусу (Naryn oblast) block -----------------------------------
$ws.Range("R22").Formula = "=R23/R24*100"
$ws.Range("R23").Value = 0.9
$ws.Range("R24").Value = 17172.7

# --- Ош облусу (Osh oblast) block ----------------------------------------
$ws.Range("R26").Formula = "=R27/R28*100"
$ws.Range("R27").Value = 15.9
$ws.Range("R28").Value = 56666.5

# --- Талас облусу (Talas oblast) block -----------------------------------
$ws.Range("R30").Formula = "=R31/R32*100"
$ws.Range("R31").Value = 58.5
$ws.Range("R32").Value = 30765.1

# --- Чүй облусу (Chui oblast) block ---------------------------------------
$ws.Range("R34").Formula = "=R35/R36*100"
$ws.Range("R35").Value = 78.3
$ws.Range("R36").Value = 110267.1

# --- Бишкек ш. (Bishkek city) block: no data available for 2021 ----------
$ws.Range("R38").Value = "-"
$ws.Range("R39").Value = "-"
$ws.Range("R40").Value = 297797.59999999998

# --- Ош ш. (Osh city) block: no data available for 2021 ------------------
$ws.Range("R42").Value = "-"
$ws.Range("R43").Value = "-"
$ws.Range("R44").Value = 39827.4

# --- Selection, matching the author's saved view -------------------------
$ws.Range("T9").Select()
